# "added new data with new plunger" - reading_vs_weight160225.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3's "3 mL" reading was a placeholder (-1000); replace it with the first
# real reading/weight pair taken with the new plunger, and fill in rows 4-5
# which previously had no "3 mL" data at all.
$ws.Range("D3").Value = 1030
$ws.Range("E3").Value = 1000.4

$ws.Range("D4").Value = 1030
$ws.Range("E4").Value = 101.5

$ws.Range("D5").Value = 1030
$ws.Range("E5").Value = 101.6

# New "1mL" reading/weight data points taken with the new plunger (rows 22-24).
$ws.Range("A22").Value = 310
$ws.Range("B22").Value = 299.10000000000002

$ws.Range("A23").Value = 310
$ws.Range("B23").Value = 300.39999999999998

$ws.Range("A24").Value = 310
$ws.Range("B24").Value = 299.5

# Leave the selection where data entry finished.
$ws.Range("A25").Select()
